# The workbook gained one new data row: a new record was inserted at row 992
# (pushing the former rows 992:1100 down to 993:1101) on the single sheet.
# Reproduce that by inserting a blank row at 992 and then populating it with
# the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 992; Excel shifts rows 992:1100 down to 993:1101 and
# extends the used range to row 1101 automatically.
$ws.Rows("992:992").Insert()

# Populate the newly inserted row 992 with the new record.
$ws.Range("A992").Value = 3
$ws.Range("B992").Value = "Femacal de La Calera"
$ws.Range("C992").Value = "Coquimbo"
$ws.Range("D992").Value = 45194
$ws.Range("E992").Value = 5
$ws.Range("F992").Value = 100112006
$ws.Range("G992").Value = "Repollo"
$ws.Range("H992").Value = "Crespo record"
$ws.Range("I992").Value = "Primera"
$ws.Range("J992").Value = 2800
$ws.Range("K992").Value = 750
$ws.Range("L992").Value = 800
$ws.Range("M992").Value = 779
$ws.Range("N992").Value = "`$/unidad"
$ws.Range("O992").Value = "Provincia de Quillota"
$ws.Range("P992").Value = 779
$ws.Range("Q992").Value = 1
$ws.Range("R992").Value = "Hortaliza"
